$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1967.7273
$ws.Range("I17").Value = 968.75
$ws.Range("J17").Value = 2538.5715
$ws.Range("K17").Value = 2906.25
$ws.Range("L17").Value = 7615.7145
$ws.Range("M17").Value = -2738.25
$ws.Range("N17").Value = -7951.7145

$ws.Range("H40").Value = 5899.9414
$ws.Range("J40").Value = 4056.2222
$ws.Range("L40").Value = 4056.2222
$ws.Range("N40").Value = -4406.2222

$ws.Range("H112").Value = 843.30554
$ws.Range("J112").Value = 815.97144
$ws.Range("L112").Value = 2447.91432
$ws.Range("N112").Value = -4663.91432

$ws.Range("H132").Value = 55236.863
$ws.Range("I132").Value = 56914.832
$ws.Range("J132").Value = 19999.5
$ws.Range("K132").Value = 170744.496
$ws.Range("L132").Value = 59998.5
$ws.Range("M132").Value = -168214.496
$ws.Range("N132").Value = -65058.5

$ws.Range("H137").Value = 4196
$ws.Range("I137").Value = 3417.8096
$ws.Range("J137").Value = 5557.8335
$ws.Range("K137").Value = 10253.4288
$ws.Range("L137").Value = 16673.5005
$ws.Range("M137").Value = -7703.4288
$ws.Range("N137").Value = -21773.5005

$ws.Range("H138").Value = 4576.727
$ws.Range("I138").Value = 7500
$ws.Range("J138").Value = 4388.129
$ws.Range("K138").Value = 22500
$ws.Range("L138").Value = 13164.387
$ws.Range("M138").Value = -17360
$ws.Range("N138").Value = -23444.387

$ws.Range("H141").Value = 977.4
$ws.Range("I141").Value = 977.4
$ws.Range("K141").Value = 2932.2
$ws.Range("M141").Value = 2247.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6495315.5
$ws.Range("I32").Value = 7937294
$ws.Range("K32").Value = 7937294
$ws.Range("M32").Value = -7937007

$ws.Range("H97").Value = 1172.7059
$ws.Range("J97").Value = 1215.1666
$ws.Range("L97").Value = 1215.1666
$ws.Range("N97").Value = -2207.1666

$ws.Range("H132").Value = 2414.83
$ws.Range("I132").Value = 2148.5745
$ws.Range("K132").Value = 6445.7235
$ws.Range("M132").Value = -3915.7235

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 37747.5
$ws.Range("J58").Value = 46995
$ws.Range("L58").Value = 46995
$ws.Range("N58").Value = -47583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 29999.75
$ws.Range("J59").Value = 29999.75
$ws.Range("L59").Value = 29999.75
$ws.Range("N59").Value = -32289.75

$ws.Range("H62").Value = 1955.4445
$ws.Range("I62").Value = 1600
$ws.Range("J62").Value = 3199.5
$ws.Range("K62").Value = 1600
$ws.Range("L62").Value = 3199.5
$ws.Range("M62").Value = -976
$ws.Range("N62").Value = -4447.5

$ws.Range("H65").Value = 1955.4445
$ws.Range("I65").Value = 1600
$ws.Range("J65").Value = 3199.5
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 15997.5
$ws.Range("M65").Value = -4880
$ws.Range("N65").Value = -22237.5

$ws.Range("H132").Value = 9260948
$ws.Range("I132").Value = 1552.4286
$ws.Range("K132").Value = 4657.2858
$ws.Range("M132").Value = -2127.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 710
$ws.Range("I8").Value = 710
$ws.Range("K8").Value = 2130
$ws.Range("M8").Value = -1991

$ws.Range("H93").Value = 5355.6665
$ws.Range("J93").Value = 6027
$ws.Range("L93").Value = 18081
$ws.Range("N93").Value = -21825

$ws.Range("H131").Value = 1454.0189
$ws.Range("J131").Value = 1475.1569
$ws.Range("L131").Value = 4425.4707
$ws.Range("N131").Value = -14505.4707

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 269482.22
$ws.Range("I80").Value = 390843.38
$ws.Range("J80").Value = 6533
$ws.Range("K80").Value = 390843.38
$ws.Range("L80").Value = 6533
$ws.Range("M80").Value = -389845.38
$ws.Range("N80").Value = -8529

$ws.Range("H83").Value = 269482.22
$ws.Range("I83").Value = 390843.38
$ws.Range("J83").Value = 6533
$ws.Range("K83").Value = 1954216.9
$ws.Range("L83").Value = 32665
$ws.Range("M83").Value = -1949224.9
$ws.Range("N83").Value = -42649

$ws.Range("H122").Value = 68174.47
$ws.Range("I122").Value = 74291.42999999999
$ws.Range("K122").Value = 222874.29
$ws.Range("M122").Value = -220424.29

$ws.Range("H132").Value = 23268942
$ws.Range("I132").Value = 34494732
$ws.Range("K132").Value = 103484196
$ws.Range("M132").Value = -103481666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H61").Value = 2404.303
$ws.Range("I61").Value = 1539.9048
$ws.Range("K61").Value = 1539.9048
$ws.Range("M61").Value = -1337.9048

$ws.Range("H68").Value = 4801.1816
$ws.Range("J68").Value = 5476.625
$ws.Range("L68").Value = 5476.625
$ws.Range("N68").Value = -6974.625

$ws.Range("H71").Value = 4801.1816
$ws.Range("J71").Value = 5476.625
$ws.Range("L71").Value = 27383.125
$ws.Range("N71").Value = -34871.125

$ws.Range("H96").Value = 150197
$ws.Range("J96").Value = 150197
$ws.Range("L96").Value = 150197
$ws.Range("N96").Value = -155689

$ws.Range("H113").Value = 2404.303
$ws.Range("I113").Value = 1539.9048
$ws.Range("K113").Value = 1539.9048
$ws.Range("M113").Value = 630.0952

$ws.Range("H132").Value = 2121.75
$ws.Range("I132").Value = 1860.7916
$ws.Range("J132").Value = 3687.5
$ws.Range("K132").Value = 5582.3748
$ws.Range("L132").Value = 11062.5
$ws.Range("M132").Value = -3052.3748
$ws.Range("N132").Value = -16122.5

$ws.Range("H136").Value = 36363.164
$ws.Range("I136").Value = 1556.6
$ws.Range("J136").Value = 108876.836
$ws.Range("K136").Value = 4669.799999999999
$ws.Range("L136").Value = 326630.508
$ws.Range("M136").Value = -2119.799999999999
$ws.Range("N136").Value = -331730.508

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 83548390
$ws.Range("J4").Value = 125008160
$ws.Range("L4").Value = 125008160
$ws.Range("N4").Value = -125008386

$ws.Range("H30").Value = 99
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H132").Value = 2293.1936
$ws.Range("I132").Value = 2082.926
$ws.Range("J132").Value = 3712.5
$ws.Range("K132").Value = 6248.778
$ws.Range("L132").Value = 11137.5
$ws.Range("M132").Value = -3718.778
$ws.Range("N132").Value = -16197.5
